# Update imputed KNN result values in the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8.347
$ws.Range("C4").Value = -11.699
$ws.Range("B6").Value = 5.764999999999999
$ws.Range("B7").Value = 5.374000000000001
$ws.Range("C9").Value = -11.149
$ws.Range("C12").Value = -10.956
$ws.Range("B16").Value = 5.018
$ws.Range("C17").Value = -13.27
$ws.Range("C18").Value = -11.875
$ws.Range("C19").Value = -11.949
$ws.Range("B20").Value = 8.475999999999999
$ws.Range("C20").Value = -12.068
$ws.Range("C26").Value = -12.562
$ws.Range("B28").Value = 5.121
$ws.Range("B29").Value = 5.326
$ws.Range("C31").Value = -12.867
$ws.Range("B32").Value = 6.943
$ws.Range("C39").Value = -12.278
$ws.Range("B40").Value = 9.221
$ws.Range("C40").Value = -12.206
$ws.Range("C41").Value = -12.02
$ws.Range("C42").Value = -12.364
$ws.Range("C43").Value = -12.395
$ws.Range("B46").Value = 5.796000000000001
$ws.Range("C47").Value = -13.123
$ws.Range("C48").Value = -11.844
$ws.Range("B51").Value = 5.399
$ws.Range("B52").Value = 5.628
$ws.Range("B57").Value = 4.856
$ws.Range("B59").Value = 4.996
$ws.Range("B62").Value = 5.494
$ws.Range("C63").Value = -10.956
$ws.Range("C64").Value = -11.073
$ws.Range("B66").Value = 4.814
$ws.Range("B73").Value = 5.884
$ws.Range("B74").Value = 9.186999999999999
$ws.Range("C76").Value = -12.288
$ws.Range("C81").Value = -13.242
$ws.Range("C89").Value = -13.521
$ws.Range("B92").Value = 4.891
$ws.Range("C94").Value = -11.769
$ws.Range("B100").Value = 6.026999999999999
